# Updated Instruction Softskill Assessment
#
# For every item row (2..56) on the single worksheet:
#   - Fill in the (previously empty) SubItemStem column E with the same
#     text that used to be the Instruction ("Inwieweit trifft diese
#     Aussage aud dich zu?").
#   - Replace the Instruction column J text: row 3 becomes a single
#     blank space, every other row becomes "Bitte klicke an.".
# Finally, move the active selection to J3 (the cell last touched).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$subItemStem = "Inwieweit trifft diese Aussage aud dich zu?"
$clickInstruction = "Bitte klicke an."
$blankInstruction = " "

for ($row = 2; $row -le 56; $row++) {
    $ws.Range("E$row").Value = $subItemStem

    if ($row -eq 3) {
        $ws.Range("J$row").Value = $blankInstruction
    } else {
        $ws.Range("J$row").Value = $clickInstruction
    }
}

$ws.Range("J3").Select() | Out-Null
